# Apply updated values to column F ("dSF") for specific rows,
# reflecting a repull/recalculation of data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 2
$ws.Range("F4").Value = -5
$ws.Range("F6").Value = -5
$ws.Range("F11").Value = 7
$ws.Range("F12").Value = 1
$ws.Range("F21").Value = -1
$ws.Range("F28").Value = -2
$ws.Range("F30").Value = -4
